$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7215909361839294
$ws.Range("B1").Value = 1.646070957183838
$ws.Range("C1").Value = 3.958370447158813
$ws.Range("D1").Value = 1.665257811546326
$ws.Range("E1").Value = 0.9116083383560181
